# Generate Report for Handoff
# Updates the handback status of e0d7219b-7313-4a9a-acad-37081b86aeef.md
# from "Handed back: in sync with en-US" to "Ready for handoff" across all
# sheets, refreshes the related timestamps, records the "stale handback"
# error detail for the two locale sheets, and widens the Error Detail
# column so the message is readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4420b93fb92f81fea5a089865377219e457b4f94/e2e/e0d7219b-7313-4a9a-acad-37081b86aeef.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7dc985ae664ebb477ddf50deb42d9015bdfb32c/e2e/e0d7219b-7313-4a9a-acad-37081b86aeef.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the e0d7219b-7313-4a9a-acad-37081b86aeef.md file
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-02 12:56:25"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the e0d7219b-7313-4a9a-acad-37081b86aeef.md file
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-02 12:56:20"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the e0d7219b-7313-4a9a-acad-37081b86aeef.md file
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-02 12:56:25"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
